$d = $word.ActiveDocument

# 1. Rename the "NumCars" field label to "NumCar"
$d.Content.Find.Execute("NumCars", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NumCar", 2)

# 2. Fill in the empty "EmpRating" row in the Employees table (Table 2,
#    the row immediately after "BalDue")
$t = $d.Tables.Item(2)
$c = $t.Cell(27, 1)
$c.Range.Text = "EmpRating"
